$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "Alchi teleport ignores fallthrough-platforms"
$ws.Range("B18").Value = "Fixed"
$ws.Range("C18").Value = "Sandro"

$ws.Range("A17").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Interior.ColorIndex = 2

$ws.Range("A18").Select()

$wb.Application.CutCopyMode = $false
